# Updated symbol list on Mon Jan 16 09:14:33 UTC 2023 with GitHub Actions
# Refreshes Price (D), Volume(1h) (E) and Hora (G) columns on the crypto
# price sheet. Values are written with a leading apostrophe so Excel keeps
# them as text (matching the sheet's existing text-typed cells) instead of
# auto-converting numeric-looking strings into Number cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'299.16"
$ws.Range("E2").Value = "'1.84%"
$ws.Range("G2").Value = "'9"
$ws.Range("D3").Value = "'31.22"
$ws.Range("E3").Value = "'-0.26%"
$ws.Range("G3").Value = "'9"
$ws.Range("D4").Value = "'5.138"
$ws.Range("E4").Value = "'0.91%"
$ws.Range("G4").Value = "'9"
$ws.Range("D5").Value = "'0.08112"
$ws.Range("E5").Value = "'9.90%"
$ws.Range("G5").Value = "'9"
$ws.Range("D6").Value = "'2.499"
$ws.Range("E6").Value = "'50.26%"
$ws.Range("G6").Value = "'9"
$ws.Range("D7").Value = "'7.823"
$ws.Range("E7").Value = "'1.61%"
$ws.Range("G7").Value = "'9"
$ws.Range("D8").Value = "'3.839"
$ws.Range("E8").Value = "'1.78%"
$ws.Range("G8").Value = "'9"
$ws.Range("D9").Value = "'0.9125"
$ws.Range("E9").Value = "'-1.57%"
$ws.Range("G9").Value = "'9"
$ws.Range("D10").Value = "'0.1710"
$ws.Range("E10").Value = "'1.47%"
$ws.Range("G10").Value = "'9"
$ws.Range("D11").Value = "'0.07318"
$ws.Range("E11").Value = "'2.51%"
$ws.Range("G11").Value = "'9"
$ws.Range("D12").Value = "'0.07968"
$ws.Range("E12").Value = "'1.07%"
$ws.Range("G12").Value = "'9"
$ws.Range("D13").Value = "'0.03033"
$ws.Range("E13").Value = "'1.03%"
$ws.Range("G13").Value = "'9"
$ws.Range("D14").Value = "'0.09946"
$ws.Range("E14").Value = "'0.42%"
$ws.Range("G14").Value = "'9"
$ws.Range("D15").Value = "'0.001497"
$ws.Range("E15").Value = "'-0.11%"
$ws.Range("G15").Value = "'9"
$ws.Range("D16").Value = "'0.006006"
$ws.Range("E16").Value = "'-4.45%"
$ws.Range("G16").Value = "'9"
$ws.Range("D17").Value = "'3.495"
$ws.Range("E17").Value = "'1.11%"
$ws.Range("G17").Value = "'9"
$ws.Range("D18").Value = "'2.240"
$ws.Range("E18").Value = "'0.61%"
$ws.Range("G18").Value = "'9"
$ws.Range("E19").Value = "'0.73%"
$ws.Range("G19").Value = "'9"
$ws.Range("D20").Value = "'0.1351"
$ws.Range("E20").Value = "'0.12%"
$ws.Range("G20").Value = "'9"
$ws.Range("D21").Value = "'4.622"
$ws.Range("E21").Value = "'0.66%"
$ws.Range("G21").Value = "'9"
$ws.Range("D22").Value = "'0.1605"
$ws.Range("E22").Value = "'3.17%"
$ws.Range("G22").Value = "'9"
$ws.Range("D23").Value = "'0.04602"
$ws.Range("E23").Value = "'-0.98%"
$ws.Range("G23").Value = "'9"
$ws.Range("D24").Value = "'0.001265"
$ws.Range("E24").Value = "'3.44%"
$ws.Range("G24").Value = "'9"
$ws.Range("E25").Value = "'0.82%"
$ws.Range("G25").Value = "'9"
$ws.Range("D26").Value = "'0.0001184"
$ws.Range("E26").Value = "'-9.23%"
$ws.Range("G26").Value = "'9"
$ws.Range("D27").Value = "'0.0003441"
$ws.Range("E27").Value = "'82.88%"
$ws.Range("G27").Value = "'9"
$ws.Range("G28").Value = "'9"
$ws.Range("G29").Value = "'9"
$ws.Range("G30").Value = "'9"
$ws.Range("G31").Value = "'9"
$ws.Range("G32").Value = "'9"
$ws.Range("G33").Value = "'9"
$ws.Range("G34").Value = "'9"
$ws.Range("G35").Value = "'9"
$ws.Range("G36").Value = "'9"
$ws.Range("G37").Value = "'9"
$ws.Range("G38").Value = "'9"
$ws.Range("D39").Value = "'0.01807"
$ws.Range("E39").Value = "'9.27%"
$ws.Range("G39").Value = "'9"
$ws.Range("D40").Value = "'0.04530"
$ws.Range("E40").Value = "'3.23%"
$ws.Range("G40").Value = "'9"
$ws.Range("D41").Value = "'0.007253"
$ws.Range("E41").Value = "'2.21%"
$ws.Range("G41").Value = "'9"
$ws.Range("D42").Value = "'0.1342"
$ws.Range("E42").Value = "'1.25%"
$ws.Range("G42").Value = "'9"
$ws.Range("D43").Value = "'0.002258"
$ws.Range("E43").Value = "'7.14%"
$ws.Range("G43").Value = "'9"
$ws.Range("D44").Value = "'0.01060"
$ws.Range("E44").Value = "'-14.26%"
$ws.Range("G44").Value = "'9"
$ws.Range("D45").Value = "'0.00006280"
$ws.Range("E45").Value = "'4.33%"
$ws.Range("G45").Value = "'9"
$ws.Range("D46").Value = "'0.00000000752"
$ws.Range("E46").Value = "'0.32%"
$ws.Range("G46").Value = "'9"
$ws.Range("D47").Value = "'0.006664"
$ws.Range("E47").Value = "'-39.61%"
$ws.Range("G47").Value = "'9"
$ws.Range("G48").Value = "'9"
$ws.Range("D49").Value = "'0.00002107"
$ws.Range("E49").Value = "'0.32%"
$ws.Range("G49").Value = "'9"
$ws.Range("D50").Value = "'0.0002007"
$ws.Range("E50").Value = "'0.39%"
$ws.Range("G50").Value = "'9"
$ws.Range("G51").Value = "'9"
